# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each changed cell we either set a new numeric value, or clear the cell entirely
# (for cells whose <c> element was removed entirely in the target workbook).

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 2468.8462
$ws.Range("I70").Value = 2507.9167
$ws.Range("K70").Value = 7523.750100000001
$ws.Range("M70").Value = -7253.750100000001
$ws.Range("H73").Value = 2468.8462
$ws.Range("I73").Value = 2507.9167
$ws.Range("K73").Value = 7523.750100000001
$ws.Range("M73").Value = -6587.750100000001
$ws.Range("H80").Value = 1035.2667
$ws.Range("I80").Value = 444.25
$ws.Range("J80").Value = 1710.7142
$ws.Range("K80").Value = 1332.75
$ws.Range("L80").Value = 5132.142599999999
$ws.Range("M80").Value = -334.75
$ws.Range("N80").Value = -7128.142599999999
$ws.Range("H83").Value = 1035.2667
$ws.Range("I83").Value = 444.25
$ws.Range("J83").Value = 1710.7142
$ws.Range("K83").Value = 3998.25
$ws.Range("L83").Value = 15396.4278
$ws.Range("M83").Value = 993.75
$ws.Range("N83").Value = -25380.4278
$ws.Range("H86").Value = 6144.846
$ws.Range("I86").Value = 5254.5713
$ws.Range("K86").Value = 5254.5713
$ws.Range("M86").Value = -4131.5713
$ws.Range("H89").Value = 6144.846
$ws.Range("I89").Value = 5254.5713
$ws.Range("K89").Value = 26272.8565
$ws.Range("M89").Value = -20656.8565
$ws.Range("H98").Value = 991.25
$ws.Range("J98").Value = 989
$ws.Range("L98").Value = 989
$ws.Range("N98").Value = -3985
$ws.Range("H106").Value = 2970
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 991.25
$ws.Range("J122").Value = 989
$ws.Range("L122").Value = 2967
$ws.Range("N122").Value = -7867
$ws.Range("H132").Value = 2158.037
$ws.Range("I132").Value = 1837
$ws.Range("J132").Value = 3281.6667
$ws.Range("K132").Value = 5511
$ws.Range("L132").Value = 9845.000100000001
$ws.Range("M132").Value = -2981
$ws.Range("N132").Value = -14905.0001

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4227.7144
$ws.Range("I32").Value = 3089.1304
$ws.Range("K32").Value = 3089.1304
$ws.Range("M32").Value = -2802.1304
$ws.Range("H45").Value = 3525.5
$ws.Range("I45").Value = 1999.6666
$ws.Range("J45").Value = 3794.7646
$ws.Range("K45").Value = 1999.6666
$ws.Range("L45").Value = 3794.7646
$ws.Range("M45").Value = -1622.6666
$ws.Range("N45").Value = -4548.7646
$ws.Range("H61").Value = 825
$ws.Range("I61").Value = 514.2857
$ws.Range("K61").Value = 514.2857
$ws.Range("M61").Value = -302.2857
$ws.Range("H63").Value = 2579
$ws.Range("I63").Value = 1973.75
$ws.Range("K63").Value = 1973.75
$ws.Range("M63").Value = -1287.75
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 2579
$ws.Range("I66").Value = 1973.75
$ws.Range("K66").Value = 9868.75
$ws.Range("M66").Value = -6436.75
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H88").Value = 2437.3
$ws.Range("I88").Value = 1174.6
$ws.Range("J88").Value = 3700
$ws.Range("K88").Value = 1174.6
$ws.Range("L88").Value = 3700
$ws.Range("M88").Value = -768.5999999999999
$ws.Range("N88").Value = -4512
$ws.Range("H91").Value = 2437.3
$ws.Range("I91").Value = 1174.6
$ws.Range("J91").Value = 3700
$ws.Range("K91").Value = 1174.6
$ws.Range("L91").Value = 3700
$ws.Range("M91").Value = 229.4000000000001
$ws.Range("N91").Value = -6508
$ws.Range("H122").Value = 2795.8333
$ws.Range("I122").Value = 1266.3334
$ws.Range("K122").Value = 3799.0002
$ws.Range("M122").Value = -1349.0002
$ws.Range("H132").Value = 807.75
$ws.Range("I132").Value = 807.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2423.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 106.75
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 825
$ws.Range("I136").Value = 514.2857
$ws.Range("K136").Value = 1542.8571
$ws.Range("M136").Value = 1007.1429

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3303.923
$ws.Range("I86").Value = 2971.8572
$ws.Range("J86").Value = 3691.3333
$ws.Range("K86").Value = 2971.8572
$ws.Range("L86").Value = 3691.3333
$ws.Range("M86").Value = -1848.8572
$ws.Range("N86").Value = -5937.3333
$ws.Range("H89").Value = 3303.923
$ws.Range("I89").Value = 2971.8572
$ws.Range("J89").Value = 3691.3333
$ws.Range("K89").Value = 14859.286
$ws.Range("L89").Value = 18456.6665
$ws.Range("M89").Value = -9243.286
$ws.Range("N89").Value = -29688.6665
$ws.Range("H94").Value = 3399.8
$ws.Range("I94").Value = 999.5
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 999.5
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -548.5
$ws.Range("N94").Value = -5902
$ws.Range("H99").Value = 1683.1428
$ws.Range("I99").Value = 1271.5
$ws.Range("K99").Value = 1271.5
$ws.Range("M99").Value = 226.5
$ws.Range("H107").Value = 875.25
$ws.Range("J107").Value = 751
$ws.Range("L107").Value = 751
$ws.Range("N107").Value = -4591

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1351.2
$ws.Range("I7").Value = 737.86664
$ws.Range("J7").Value = 2271.2
$ws.Range("K7").Value = 737.86664
$ws.Range("L7").Value = 2271.2
$ws.Range("M7").Value = -624.86664
$ws.Range("N7").Value = -2497.2
$ws.Range("H88").Value = 29750
$ws.Range("J88").Value = 29750
$ws.Range("L88").Value = 29750
$ws.Range("N88").Value = -30562
$ws.Range("H91").Value = 29750
$ws.Range("J91").Value = 29750
$ws.Range("L91").Value = 29750
$ws.Range("N91").Value = -32558

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 157.8
$ws.Range("I38").Value = 146.33333
$ws.Range("K38").Value = 438.99999
$ws.Range("M38").Value = -91.99998999999997

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3360.111
$ws.Range("I102").Value = 3342.625
$ws.Range("K102").Value = 3342.625
$ws.Range("M102").Value = -1720.625

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 557
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 699.5
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 699.5
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1289.5
$ws.Range("H27").Value = 557
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 699.5
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 699.5
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -913.5
$ws.Range("H46").Value = 4103.8696
$ws.Range("J46").Value = 3949.0833
$ws.Range("L46").Value = 3949.0833
$ws.Range("N46").Value = -4325.0833
$ws.Range("H82").Value = 840.2
$ws.Range("I82").Value = 880.6
$ws.Range("J82").Value = 799.8
$ws.Range("K82").Value = 880.6
$ws.Range("L82").Value = 799.8
$ws.Range("M82").Value = -519.6
$ws.Range("N82").Value = -1521.8
$ws.Range("H85").Value = 840.2
$ws.Range("I85").Value = 880.6
$ws.Range("J85").Value = 799.8
$ws.Range("K85").Value = 880.6
$ws.Range("L85").Value = 799.8
$ws.Range("M85").Value = 367.4
$ws.Range("N85").Value = -3295.8
$ws.Range("H136").Value = 2252.35
$ws.Range("I136").Value = 2254.0625
$ws.Range("J136").Value = 2245.5
$ws.Range("K136").Value = 6762.1875
$ws.Range("L136").Value = 6736.5
$ws.Range("M136").Value = -4212.1875
$ws.Range("N136").Value = -11836.5

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1806.125
$ws.Range("J96").Value = 1683.1666
$ws.Range("L96").Value = 1683.1666
$ws.Range("N96").Value = -4429.1666
$ws.Range("H100").Value = 2063
$ws.Range("I100").Value = 1339.6666
$ws.Range("J100").Value = 3799
$ws.Range("K100").Value = 2679.3332
$ws.Range("L100").Value = 7598
$ws.Range("M100").Value = -2138.3332
$ws.Range("N100").Value = -8680
$ws.Range("H122").Value = 3006.5715
$ws.Range("I122").Value = 2637.25
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 7911.75
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -5461.75

